$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "ZeroMQ for messaging " -> "ZMQ for messaging " in the
# Test Environment bullet list, also dropping the (now stale)
# spell-check proofErr markers that wrapped "ZeroMQ".
# ------------------------------------------------------------------
$zmqPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "ZeroMQ for messaging*") {
        $zmqPara = $cand
        break
    }
}

if ($zmqPara -ne $null) {
    # Insert a throw-away marker character immediately before the run so
    # the subsequent Find/Replace range no longer starts exactly on top of
    # the orphaned <w:proofErr w:type="spellStart"/> marker (which,
    # otherwise, survives any replace that begins precisely at that
    # boundary). Swallowing it inside a wider match removes it cleanly,
    # along with its matching spellEnd.
    $startPos = $zmqPara.Range.Start
    $marker = $d.Range($startPos, $startPos)
    $marker.InsertBefore("#")

    $zmqPara = $d.Paragraphs.Item($zmqPara.Index)
    $find = $zmqPara.Range.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute(
        "#ZeroMQ for messaging ", $false, $false, $false, $false, $false,
        $true, 1, $false, "ZMQ for messaging ", 2
    )
}

# ------------------------------------------------------------------
# Change 2: the blank spacer paragraph right before "7. Code-Based
# Testing Strategy" moves from the Heading2 style to a List Bullet
# paragraph with its numbering switched off (numId=0) and a manual
# 720-twip left indent.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $candText = $cand.Range.Text.TrimEnd([char]13)
    if ($candText -eq "" -and $cand.Style.NameLocal -eq "Heading 2") {
        $next = $cand.Next()
        if ($next -ne $null -and $next.Range.Text -like "7. Code-Based Testing Strategy*") {
            $target = $cand
            break
        }
    }
}

if ($target -ne $null) {
    $target.Style = "List Bullet"
    $target.Range.ListFormat.RemoveNumbers()
    $target.LeftIndent = 36
}
